# ---------------------------------------------------------------------------
# "new results plus statistics"
#
# Adds two new result blocks (RF / ST) to the right of the existing
# Decision-Tree (C:G/D6:G13-ish) and Random-Forest (K:O) comparison tables,
# each with Global / Local / GA sub-columns, plus marks a couple of existing
# MAE cells (G14, M14) bold to highlight the winning value.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Two existing "osales" MAE cells get bolded (best-value highlight).
# ---------------------------------------------------------------------------
$ws.Range("G14").Font.Bold = $true
$ws.Range("M14").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) New merged group headers in row 5: "RF" over S5:U5, "ST" over Y5:AA5.
# ---------------------------------------------------------------------------
$ws.Range("S5").Value = "RF"
$ws.Range("Y5").Value = "ST"
$ws.Range("S5:U5").Merge()
$ws.Range("Y5:AA5").Merge()
$ws.Range("S5:U5").HorizontalAlignment = -4108
$ws.Range("Y5:AA5").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 3) Row 6 sub-headers (Global / Local / GA) for each new block.
# ---------------------------------------------------------------------------
$ws.Range("S6").Value = "Global"
$ws.Range("T6").Value = "Local"
$ws.Range("U6").Value = "GA"
$ws.Range("Y6").Value = "Global"
$ws.Range("Z6").Value = "Local"
$ws.Range("AA6").Value = "GA"

# ---------------------------------------------------------------------------
# 4) Data rows 7-14 for the "RF" block (S:U) and "ST" block (Y:AA).
#    Column layout within each block: Global, Local, GA.
#    Bold marks the "winning"/highlighted value in each row, mirroring the
#    existing D:G / L:O tables.
# ---------------------------------------------------------------------------

function Set-Stat($cell, $value, $bold) {
    $r = $ws.Range($cell)
    $r.Value = $value
    $r.NumberFormat = "0.000"
    if ($bold) {
        $r.Font.Bold = $true
    }
}

# Row 7
Set-Stat "S7" 240.69377499999999 $false
Set-Stat "T7" 227.76239749999999 $true
Set-Stat "U7" 229.98624375 $false
Set-Stat "Y7" 335.52234437499999 $false
Set-Stat "Z7" 335.53621562500001 $false
Set-Stat "AA7" 293.04720424999999 $true

# Row 8
Set-Stat "S8" 395.10045687500002 $false
Set-Stat "T8" 379.34393125000003 $true
Set-Stat "U8" 385.13324875000001 $false
Set-Stat "Y8" 507.864773125 $false
Set-Stat "Z8" 480.957105625 $false
Set-Stat "AA8" 477.02404349999898 $true

# Row 9
Set-Stat "S9" 0.42858000000000002 $false
Set-Stat "T9" 0.41279500000000002 $true
Set-Stat "U9" 0.417603749999999 $false
Set-Stat "Y9" 5.0400912499999997 $false
Set-Stat "Z9" 0.60032874999999997 $true
Set-Stat "AA9" 0.60081374999999904 $false

# Row 10
Set-Stat "S10" 0.87112999999999996 $false
Set-Stat "T10" 0.72140249999999995 $true
Set-Stat "U10" 0.77481250000000002 $false
Set-Stat "Y10" 5.0432100000000002 $false
Set-Stat "Z10" 0.61703624999999995 $false
Set-Stat "AA10" 0.61270000000000002 $true

# Row 11
Set-Stat "S11" 53.587011250000003 $false
Set-Stat "T11" 49.55937625 $true
Set-Stat "U11" 49.738267499999999 $false
Set-Stat "Y11" 79.866088750000003 $false
Set-Stat "Z11" 68.365613124999996 $false
Set-Stat "AA11" 68.187651499999902 $true

# Row 12
Set-Stat "S12" 66.070739375000002 $false
Set-Stat "T12" 64.811838750000007 $true
Set-Stat "U12" 64.859575625000005 $false
Set-Stat "Y12" 93.221359375000006 $false
Set-Stat "Z12" 85.271510000000006 $false
Set-Stat "AA12" 84.200227124999998 $true

# Row 13
Set-Stat "S13" 0.83244571428571401 $false
Set-Stat "T13" 0.82558428571428599 $false
Set-Stat "U13" 0.82457499999999995 $true
Set-Stat "Y13" 0.90133357142857096 $false
Set-Stat "Z13" 0.86115928571428602 $true
Set-Stat "AA13" 0.86744185714285704 $false

# Row 14
Set-Stat "S14" 2959.8802258333299 $false
Set-Stat "T14" 2938.1286341666701 $true
Set-Stat "U14" 2945.4223541666602 $false
Set-Stat "Y14" 2987.2809833333299 $false
Set-Stat "Z14" 3270.1264875000002 $false
Set-Stat "AA14" 2945.4223541666602 $true

# ---------------------------------------------------------------------------
# 5) New column widths for L:M (best-fit-style custom width) and selection
#    on the new "RF" header, scrolled into view.
# ---------------------------------------------------------------------------
$ws.Range("L1:M1").ColumnWidth = 11.33

$ws.Range("S5:U5").Select()
